$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''20.226.75'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '''1.440.97'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("D4").Value = '''1.008'
$ws.Range("E4").Value = '  +0.89%  '
$ws.Range("D5").Value = '''0.9212'
$ws.Range("E5").Value = '  -7.90%  '
$ws.Range("D6").Value = '''274.83'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '''0.3636'
$ws.Range("E7").Value = '  -1.49%  '
$ws.Range("D8").Value = '''0.3068'
$ws.Range("E8").Value = '  -1.89%  '
$ws.Range("D9").Value = '''38.65'
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").Value = '''1.014'
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").Value = '''0.06454'
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = '''0.9993'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '''5.294'
$ws.Range("E13").Value = '  -3.18%  '
$ws.Range("D14").Value = '''17.31'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '''6.002'
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("D16").Value = '''0.00001002'
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").Value = '''1.442.40'
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("D18").Value = '''0.9390'
$ws.Range("E18").Value = '  -6.14%  '
$ws.Range("D19").Value = '''0.05613'
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("D20").Value = '''67.46'
$ws.Range("E20").Value = '  -4.58%  '
$ws.Range("D21").Value = '''5.305'
$ws.Range("E21").Value = '  -5.38%  '
$ws.Range("D22").Value = '''14.12'
$ws.Range("E22").Value = '  -4.36%  '
$ws.Range("D23").Value = '''10.69'
$ws.Range("E23").Value = '  -3.44%  '
$ws.Range("D24").Value = '''2.240'
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").Value = '''20.261.35'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").Value = '''138.70'
$ws.Range("E26").Value = '  +2.52%  '
$ws.Range("D27").Value = '''2.030'
$ws.Range("E27").Value = '  -10.56%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = '''1.594.82'
$ws.Range("E29").Value = '  +1.86%  '
$ws.Range("D30").Value = '''109.94'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("D32").Value = '''4.775'
$ws.Range("E32").Value = '  -10.89%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '''0.07648'
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7693'
$ws.Range("E34").Value = '  -6.93%  '
$ws.Range("D35").Value = '''1.453'
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("D36").Value = '''0.05720'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("D37").Value = '''1.132'
$ws.Range("E37").Value = '  +3.86%  '
$ws.Range("D38").Value = '''4.604'
$ws.Range("E38").Value = '  -5.82%  '
$ws.Range("D39").Value = '''0.01974'
$ws.Range("E39").Value = '  -4.69%  '
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = '''0.9294'
$ws.Range("E40").Value = '  -7.07%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '''10.08'
$ws.Range("E41").Value = '  -4.13%  '
$ws.Range("D42").Value = '''0.1827'
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D43").Value = '''6.920'
$ws.Range("E43").Value = '  -17.66%  '
$ws.Range("B44").Value = 'PancakeSwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D44").Value = '''3.474'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.5166'
$ws.Range("E45").Value = '  -2.75%  '
$ws.Range("D46").Value = '''11.67'
$ws.Range("E46").Value = '  -4.85%  '
$ws.Range("D47").Value = '''114.69'
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("D48").Value = '''0.5026'
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("D49").Value = '''1.717'
$ws.Range("E49").Value = '  -3.24%  '
$ws.Range("D50").Value = '''0.06350'
$ws.Range("E50").Value = '  +2.55%  '
$ws.Range("D51").Value = '''0.9863'
$ws.Range("E51").Value = '  -1.31%  '
